# Add a new "2022-Q3" sheet (holdings detail) right after the "总计" sheet,
# and insert a corresponding summary row into the "总计" sheet.

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)

# --- 1. Create the new "2022-Q3" worksheet right after "总计" ---
$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q3"

# Re-fetch the previously-first quarter sheet (now shifted) by name so we
# have a fresh, valid reference to copy formatting from.
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# Copy header row (values + style) from the 2022-Q2 sheet.
$q2Sheet.Range("B1:H1").Copy($newSheet.Range("B1"))

# Copy row-2 formatting only (keeps the bordered/bold "index" style on A2
# and default style on the rest) before writing the real Q3 values in.
$q2Sheet.Range("A2:H2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Fill in the 2022-Q3 holder detail row.
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "008116"
$newSheet.Range("C2").Value = "银华沪深股通精选混合"
$newSheet.Range("D2").Value = "'0.50"
$newSheet.Range("E2").Value = "'88.46"
$newSheet.Range("F2").Value = "'4.04"
$newSheet.Range("G2").Value = "'0.0202"
$newSheet.Range("H2").Value = 8

# --- 2. Insert the new "2022-Q3" summary row into the "总计" sheet ---
$totalSheet.Rows.Item(2).Insert()

# Excel's row-insert carries over formatting from the row above; clear it
# so the new row starts unstyled, matching the rest of the data rows.
$totalSheet.Rows.Item(2).ClearFormats()

# Give the new A2 the same "index" style used by the other rows (copy
# format only from A3, which already carries the style after the shift).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.02

# Renumber the index column (A) for the rows that shifted down, so it
# stays a contiguous 0-based sequence.
for ($r = 3; $r -le 9; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

Write-Output "2022-Q3 sheet added and 总计 updated"
